$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 2632727
$ws.Range("I129").Value = 31251008
$ws.Range("J129").Value = 1160.816
$ws.Range("K129").Value = 93753024
$ws.Range("L129").Value = 3482.448
$ws.Range("M129").Value = -93748024
$ws.Range("N129").Value = -13482.448

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1557.3422
$ws.Range("I45").Value = 999.61536
$ws.Range("J45").Value = 2765.75
$ws.Range("K45").Value = 999.61536
$ws.Range("L45").Value = 2765.75
$ws.Range("M45").Value = -622.61536
$ws.Range("N45").Value = -3519.75

$ws.Range("H74").Value = 679.34375
$ws.Range("I74").Value = 643.7586
$ws.Range("J74").Value = 1023.3333
$ws.Range("K74").Value = 643.7586
$ws.Range("L74").Value = 1023.3333
$ws.Range("M74").Value = 230.2414
$ws.Range("N74").Value = -2771.3333

$ws.Range("H77").Value = 679.34375
$ws.Range("I77").Value = 643.7586
$ws.Range("J77").Value = 1023.3333
$ws.Range("K77").Value = 3218.793
$ws.Range("L77").Value = 5116.6665
$ws.Range("M77").Value = 1149.207
$ws.Range("N77").Value = -13852.6665

$ws.Range("H97").Value = 693.4737
$ws.Range("I97").Value = 661
$ws.Range("J97").Value = 866.6667
$ws.Range("K97").Value = 661
$ws.Range("L97").Value = 866.6667
$ws.Range("M97").Value = -165
$ws.Range("N97").Value = -1858.6667

$ws.Range("H132").Value = 20002600
$ws.Range("I132").Value = 29413494
$ws.Range("J132").Value = 4448.625
$ws.Range("K132").Value = 88240482
$ws.Range("L132").Value = 13345.875
$ws.Range("M132").Value = -88237952
$ws.Range("N132").Value = -18405.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 35000
$ws.Range("J59").Value = 35000
$ws.Range("L59").Value = 35000
$ws.Range("N59").Value = -36694

$ws.Range("H107").Value = 4043.25
$ws.Range("I107").Value = 2277.6667
$ws.Range("J107").Value = 5102.6
$ws.Range("K107").Value = 2277.6667
$ws.Range("L107").Value = 5102.6
$ws.Range("M107").Value = -357.6667000000002
$ws.Range("N107").Value = -8942.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 13515930
$ws.Range("I58").Value = 969.75
$ws.Range("J58").Value = 38466624
$ws.Range("K58").Value = 969.75
$ws.Range("L58").Value = 38466624
$ws.Range("M58").Value = -766.75
$ws.Range("N58").Value = -38467030

$ws.Range("H70").Value = 34249.75
$ws.Range("J70").Value = 34000
$ws.Range("L70").Value = 34000
$ws.Range("N70").Value = -34630

$ws.Range("H73").Value = 34249.75
$ws.Range("J73").Value = 34000
$ws.Range("L73").Value = 34000
$ws.Range("N73").Value = -36184

$ws.Range("H105").Value = 2688.875
$ws.Range("I105").Value = 2130.1428
$ws.Range("J105").Value = 3123.4443
$ws.Range("K105").Value = 2130.1428
$ws.Range("L105").Value = 3123.4443
$ws.Range("M105").Value = -383.1428000000001
$ws.Range("N105").Value = -6617.4443

$ws.Range("H132").Value = 2018.1842
$ws.Range("I132").Value = 1543.4445
$ws.Range("J132").Value = 3183.4546
$ws.Range("K132").Value = 4630.333500000001
$ws.Range("L132").Value = 9550.363799999999
$ws.Range("M132").Value = -2100.333500000001
$ws.Range("N132").Value = -14610.3638

$ws.Range("H134").Value = 1564.5927
$ws.Range("I134").Value = 918.2105
$ws.Range("J134").Value = 3099.75
$ws.Range("K134").Value = 2754.6315
$ws.Range("L134").Value = 9299.25
$ws.Range("M134").Value = -219.6315
$ws.Range("N134").Value = -14369.25

$ws.Range("H136").Value = 13515930
$ws.Range("I136").Value = 969.75
$ws.Range("J136").Value = 38466624
$ws.Range("K136").Value = 2909.25
$ws.Range("L136").Value = 115399872
$ws.Range("M136").Value = -359.25
$ws.Range("N136").Value = -115404972

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1117.85
$ws.Range("I5").Value = 644.5294
$ws.Range("J5").Value = 3800
$ws.Range("K5").Value = 1933.5882
$ws.Range("L5").Value = 11400
$ws.Range("M5").Value = -1821.5882
$ws.Range("N5").Value = -11624

$ws.Range("H12").Value = 133.16
$ws.Range("I12").Value = 7.8
$ws.Range("J12").Value = 164.5
$ws.Range("K12").Value = 23.4
$ws.Range("L12").Value = 493.5
$ws.Range("M12").Value = 149.6
$ws.Range("N12").Value = -839.5

$ws.Range("H122").Value = 1828.3846
$ws.Range("I122").Value = 635
$ws.Range("J122").Value = 2574.25
$ws.Range("K122").Value = 5715
$ws.Range("L122").Value = 23168.25
$ws.Range("M122").Value = -3265
$ws.Range("N122").Value = -28068.25

$ws.Range("H131").Value = 1597
$ws.Range("I131").Value = 6294
$ws.Range("J131").Value = 1149.6666
$ws.Range("K131").Value = 18882
$ws.Range("L131").Value = 3448.9998
$ws.Range("M131").Value = -13842
$ws.Range("N131").Value = -13528.9998

$ws.Range("H135").Value = 1117.85
$ws.Range("I135").Value = 644.5294
$ws.Range("J135").Value = 3800
$ws.Range("K135").Value = 5800.7646
$ws.Range("L135").Value = 34200
$ws.Range("M135").Value = -3265.7646
$ws.Range("N135").Value = -39270

$ws.Range("H136").Value = 2025.4348
$ws.Range("I136").Value = 1333.0625
$ws.Range("J136").Value = 3608
$ws.Range("K136").Value = 3999.1875
$ws.Range("L136").Value = 10824
$ws.Range("M136").Value = 1100.8125
$ws.Range("N136").Value = -21024

$ws.Range("H137").Value = 2535.1
$ws.Range("I137").Value = 1706.3636
$ws.Range("J137").Value = 3548
$ws.Range("K137").Value = 5119.0908
$ws.Range("L137").Value = 10644
$ws.Range("M137").Value = -19.09079999999994
$ws.Range("N137").Value = -20844

$ws.Range("H139").Value = 8338967
$ws.Range("I139").Value = 14707842
$ws.Range("J139").Value = 10438.615
$ws.Range("K139").Value = 44123526
$ws.Range("L139").Value = 31315.845
$ws.Range("M139").Value = -44118386
$ws.Range("N139").Value = -41595.845

$ws.Range("H140").Value = 8335880
$ws.Range("I140").Value = 15152375
$ws.Range("J140").Value = 4607.778
$ws.Range("K140").Value = 45457125
$ws.Range("L140").Value = 13823.334
$ws.Range("M140").Value = -45451945
$ws.Range("N140").Value = -24183.334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 30000
$ws.Range("J34").Value = 30000
$ws.Range("L34").Value = 30000
$ws.Range("N34").Value = -30536

$ws.Range("H76").Value = 30000
$ws.Range("J76").Value = 30000
$ws.Range("L76").Value = 30000
$ws.Range("N76").Value = -30630

$ws.Range("H79").Value = 30000
$ws.Range("J79").Value = 30000
$ws.Range("L79").Value = 30000
$ws.Range("N79").Value = -32184

$ws.Range("H102").Value = 55401.684
$ws.Range("I102").Value = 2289.5625
$ws.Range("K102").Value = 2289.5625
$ws.Range("M102").Value = -667.5625

$ws.Range("H122").Value = 5498.1577
$ws.Range("I122").Value = 4384.1113
$ws.Range("J122").Value = 6500.8
$ws.Range("K122").Value = 13152.3339
$ws.Range("L122").Value = 19502.4
$ws.Range("M122").Value = -10702.3339
$ws.Range("N122").Value = -24402.4

$ws.Range("H126").Value = 4164.909
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 4646
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 13938
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -18878

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2573.6365
$ws.Range("I7").Value = 1761.4286
$ws.Range("J7").Value = 3995
$ws.Range("K7").Value = 1761.4286
$ws.Range("L7").Value = 3995
$ws.Range("M7").Value = -1649.4286
$ws.Range("N7").Value = -4219

$ws.Range("H126").Value = 2573.6365
$ws.Range("I126").Value = 1761.4286
$ws.Range("J126").Value = 3995
$ws.Range("K126").Value = 5284.2858
$ws.Range("L126").Value = 11985
$ws.Range("M126").Value = -2814.2858
$ws.Range("N126").Value = -16925

$ws.Range("H137").Value = 29495
$ws.Range("J137").Value = 29495
$ws.Range("L137").Value = 29495
$ws.Range("N137").Value = -39695

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 171147.72
$ws.Range("J21").Value = 66011.336
$ws.Range("L21").Value = 66011.336
$ws.Range("N21").Value = -66481.336

$ws.Range("H35").Value = 171147.72
$ws.Range("J35").Value = 66011.336
$ws.Range("L35").Value = 66011.336
$ws.Range("N35").Value = -66591.336

$ws.Range("H126").Value = 4350730.5
$ws.Range("I126").Value = 2164.5881
$ws.Range("K126").Value = 6493.7643
$ws.Range("M126").Value = -4023.7643

$ws.Range("H132").Value = 4643.6045
$ws.Range("I132").Value = 1756.1428
$ws.Range("J132").Value = 10033.533
$ws.Range("K132").Value = 5268.428400000001
$ws.Range("L132").Value = 30100.599
$ws.Range("M132").Value = -2738.428400000001
$ws.Range("N132").Value = -35160.599
